$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.396.83'

$ws.Range("D3").Value = '1.596.75'

$ws.Range("E4").Value = '  -0.24%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.77'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.93%  '

$ws.Range("E6").Value = '  +1.35%  '

$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.09'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +8.75%  '

$ws.Range("E9").Value = '  +0.74%  '

$ws.Range("E10").Value = '  +0.99%  '

$ws.Range("E11").Value = '  +1.89%  '

$ws.Range("D12").Value = '1.823.51'
$ws.Range("E12").Value = '  +2.00%  '

$ws.Range("D13").Value = '1.595.91'
$ws.Range("E13").Value = '  +2.00%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.76'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.02%  '

$ws.Range("D16").Value = '28.406.23'
$ws.Range("E16").Value = '  +4.41%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.24'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.11%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '228.66'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.91%  '

$ws.Range("E19").Value = '  +1.56%  '

$ws.Range("E21").Value = '  -0.13%  '

$ws.Range("E22").Value = '  -0.72%  '

$ws.Range("E23").Value = '  -0.23%  '

$ws.Range("E24").Value = '  +0.58%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.27%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.23'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.34%  '

$ws.Range("E27").Value = '  +0.71%  '

$ws.Range("E28").Value = '  -0.44%  '

$ws.Range("E29").Value = '  -0.06%  '

$ws.Range("E30").Value = '  +1.07%  '

$ws.Range("E31").Value = '  +1.22%  '

$ws.Range("E32").Value = '  +0.03%  '

$ws.Range("E33").Value = '  -0.39%  '

$ws.Range("D34").Value = '1.399.57'
$ws.Range("E34").Value = '  -4.02%  '

$ws.Range("E35").Value = '  -1.24%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.05'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.11%  '

$ws.Range("E37").Value = '  +0.35%  '

$ws.Range("E38").Value = '  +0.89%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.53'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.24%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.541'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.23%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.817'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.42%  '

$ws.Range("E42").Value = '  -2.47%  '

$ws.Range("E43").Value = '  -0.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.90'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.51%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.982'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.38%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.54'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.19%  '

$ws.Range("D47").Value = '1.733.72'
$ws.Range("E47").Value = '  +2.02%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.69'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.13%  '

$ws.Range("E49").Value = '  -0.15%  '

$ws.Range("E50").Value = '  -1.16%  '

$ws.Range("E51").Value = '  +0.11%  '

